$d = $word.ActiveDocument

$replacements = @(
    @("0.430031698", "0.4300223119"),
    @("0.001000320", "0.0009996044"),
    @("0.867003096", "0.8670618968"),
    @("0.004996671", "0.0049992289"),
    @("0.901021423", "0.9009734294"),
    @("0.979989177", "0.9799282990"),
    @("0.800105576", "0.7999964995"),
    @("0.509920051", "0.5099962315"),
    @("0.879920527", "0.8802751683")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
